$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the shared formulas in column A and set plain values
$ws.Range("A1").Value = 5
$ws.Range("A2").Value = 15
$ws.Range("A3").Value = 25
$ws.Range("A4").Value = 35

# Update column B values
$ws.Range("B1").Value = 3
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 2

# Update selection to B5 (single cell)
$ws.Range("B5").Select()
